$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Dave Cloete"

# 2. Header row 2: relabel E2, add F2
$ws.Range("E2").Value = "Estimated Hours"
$ws.Range("F2").Value = "Actual Hours"

# 3. Data rows 3-29: task names (A), estimated hours (E); row 29 is brand new so needs B/C/D too
$ws.Range("A3").Value = "Project Group Meeting 1"
$ws.Range("E3").Value = 2
$ws.Range("A4").Value = "Setup PI with Node.js and NAS mount"
$ws.Range("E4").Value = 2
$ws.Range("A5").Value = "UI POC: Node.js Static web and Controllers"
$ws.Range("E5").Value = 2
$ws.Range("A6").Value = "Project Group Meeting 2"
$ws.Range("E6").Value = 2
$ws.Range("A7").Value = "UI POC: Node.js Controllers - Added twig templating"
$ws.Range("E7").Value = 4
$ws.Range("A8").Value = "UI POC: D3 Charting - Globe Chart"
$ws.Range("E8").Value = 4
$ws.Range("A9").Value = "Project Group Meeting 3"
$ws.Range("E9").Value = 2
$ws.Range("A10").Value = "UI POC: D3 charting - Bar Chart tutorial, customised V Bar chart"
$ws.Range("E10").Value = 4
$ws.Range("A11").Value = "Project Group Meeting 4"
$ws.Range("E11").Value = 2
$ws.Range("A12").Value = "UI POC: D3 charting - Heat Maps (Days/Category)"
$ws.Range("E12").Value = 4
$ws.Range("A13").Value = "UI POC: D3 charting - Heat Maps (Days/months) + refactoring of "
$ws.Range("E13").Value = 4
$ws.Range("A14").Value = "Project Group Meeting 5"
$ws.Range("E14").Value = 2
$ws.Range("A15").Value = "UI POC: Refactor and Testing"
$ws.Range("E15").Value = 2
$ws.Range("A16").Value = "UI: Heat maps - Final"
$ws.Range("E16").Value = 8
$ws.Range("A17").Value = "UI: World /state Maps - Part 1"
$ws.Range("E17").Value = 2
$ws.Range("A18").Value = "Project Group Meeting 6"
$ws.Range("E18").Value = 2
$ws.Range("A19").Value = "UI: Criteria Configuration - Part 1"
$ws.Range("E19").Value = 2
$ws.Range("A20").Value = "Project Group Meeting 7"
$ws.Range("E20").Value = 2
$ws.Range("A21").Value = "UI: Criteria Configuration - Part 2"
$ws.Range("E21").Value = 4
$ws.Range("A22").Value = "UI: Initial Stream Graph Research"
$ws.Range("E22").Value = 4
$ws.Range("A23").Value = "UI: Stream Graphs - Part 1"
$ws.Range("E23").Value = 4
$ws.Range("A24").Value = "UI: Stream Graphs - Part 2"
$ws.Range("E24").Value = 4
$ws.Range("A25").Value = "Project Group Meeting 8"
$ws.Range("E25").Value = 8
$ws.Range("A26").Value = "Project Group Meeting 9"
$ws.Range("E26").Value = 2
$ws.Range("A27").Value = "Project Group Meeting 10"
$ws.Range("E27").Value = 2
$ws.Range("A28").Value = "Reports - Individual and Group & Group Meeting 11(Virtual)"
$ws.Range("E28").Value = 8
$ws.Range("A29").Value = "Project Group Meeting 12"
$ws.Range("B29").Value = 42552
$ws.Range("C29").Value = 0.4166666666666667
$ws.Range("D29").Value = 0.7083333333333334
$ws.Range("E29").Value = 8

# 4. Column F: Actual Hours = D - C, for every data row through row 44 (shared formula range like old E)
$ws.Range("F3:F44").Formula = "=D3-C3"

# 5. Sum row 45
$ws.Range("E45").Formula = "=SUM(E3:E44)"
$ws.Range("F45").Formula = "=SUM(F3:F44)"

# 6. Merge the title row across the new column range
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1:F1").Merge()

# 7. Column widths + number formats + centering for E and F
$ws.Columns.Item(5).ColumnWidth = 15.7109375
$ws.Range("E3:E44").NumberFormat = "0"
$ws.Range("E3:E44").HorizontalAlignment = -4108
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E45").HorizontalAlignment = -4108

$ws.Columns.Item(6).ColumnWidth = 14.5703125
$ws.Range("F3:F45").NumberFormat = "h:mm:ss"
$ws.Range("F3:F45").HorizontalAlignment = -4108
$ws.Range("F2").HorizontalAlignment = -4108

# 8. Selection matches author workflow (cursor left on sum cell, view scrolled back to top)
$ws.Range("E45").Select()
